# CheckboxLink/sample.xlsx -- check four of the form-control checkboxes.
#
# Sheet1 has eight Forms checkboxes (Check Box 1..8), each linked via its
# ControlFormat.LinkedCell to a cell on the sheet:
#   Check Box 1 -> B33   Check Box 4 -> C35
#   Check Box 2 -> C34   Check Box 5 -> D35
#   Check Box 3 -> B35   Check Box 6 -> C36
#                        Check Box 7 -> D36
#                        Check Box 8 -> G36
#
# Checking a linked checkbox writes TRUE into its linked cell and flips the
# control itself to the checked state. Here boxes 1, 4, 5 and 6 get checked,
# i.e. B33, C35, D35 and C36 become TRUE.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$checkedBoxes = @{
    "Check Box 1" = "B33"
    "Check Box 4" = "C35"
    "Check Box 5" = "D35"
    "Check Box 6" = "C36"
}

foreach ($name in $checkedBoxes.Keys) {
    $cellAddress = $checkedBoxes[$name]

    # Flip the control to "checked" ...
    $shape = $ws.Shapes($name)
    $shape.ControlFormat.Value = 1

    # ... and mirror that into the cell it's linked to, same as Excel does
    # for a Forms checkbox with a LinkedCell.
    $ws.Range($cellAddress).Value = $true
}
